$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1181
$ws.Range("I62").Value = 950
$ws.Range("K62").Value = 950
$ws.Range("M62").Value = -326

$ws.Range("H65").Value = 1181
$ws.Range("I65").Value = 950
$ws.Range("K65").Value = 4750
$ws.Range("M65").Value = -1630

$ws.Range("H136").Value = 54795
$ws.Range("J136").Value = 54795
$ws.Range("L136").Value = 54795
$ws.Range("N136").Value = -64995

$ws.Range("H137").Value = 5883070.5
$ws.Range("I137").Value = 682.14813
$ws.Range("J137").Value = 28572282
$ws.Range("K137").Value = 2046.44439
$ws.Range("L137").Value = 85716846
$ws.Range("M137").Value = 503.5556099999999
$ws.Range("N137").Value = -85721946

$ws.Range("H138").Value = 1601.766
$ws.Range("I138").Value = 1387.5428
$ws.Range("J138").Value = 2226.5833
$ws.Range("K138").Value = 4162.6284
$ws.Range("L138").Value = 6679.749899999999
$ws.Range("M138").Value = 977.3716000000004
$ws.Range("N138").Value = -16959.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1105.1538
$ws.Range("I2").Value = 1116.875
$ws.Range("K2").Value = 1116.875
$ws.Range("M2").Value = -1003.875

$ws.Range("H32").Value = 7118.6665
$ws.Range("I32").Value = 6479.25
$ws.Range("J32").Value = 12234
$ws.Range("K32").Value = 6479.25
$ws.Range("L32").Value = 12234
$ws.Range("M32").Value = -6192.25
$ws.Range("N32").Value = -12808

$ws.Range("H37").Value = 8007.3335
$ws.Range("I37").Value = 7034
$ws.Range("J37").Value = 8202
$ws.Range("K37").Value = 7034
$ws.Range("L37").Value = 8202
$ws.Range("M37").Value = -6761
$ws.Range("N37").Value = -8748

$ws.Range("H116").Value = 1105.1538
$ws.Range("I116").Value = 1116.875
$ws.Range("K116").Value = 1116.875
$ws.Range("M116").Value = 1177.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1105.1538
$ws.Range("I3").Value = 1116.875
$ws.Range("K3").Value = 1116.875
$ws.Range("M3").Value = -1002.875

$ws.Range("H20").Value = 4562.5
$ws.Range("I20").Value = 4416.6665
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 4416.6665
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -4169.6665
$ws.Range("N20").Value = -5494

$ws.Range("H22").Value = 419.16666
$ws.Range("I22").Value = 421
$ws.Range("K22").Value = 421
$ws.Range("M22").Value = -248

$ws.Range("H86").Value = 22729352
$ws.Range("I86").Value = 1776.1428
$ws.Range("J86").Value = 62502612
$ws.Range("K86").Value = 1776.1428
$ws.Range("L86").Value = 62502612
$ws.Range("M86").Value = -653.1428000000001
$ws.Range("N86").Value = -62504858

$ws.Range("H89").Value = 22729352
$ws.Range("I89").Value = 1776.1428
$ws.Range("J89").Value = 62502612
$ws.Range("K89").Value = 8880.714
$ws.Range("L89").Value = 312513060
$ws.Range("M89").Value = -3264.714
$ws.Range("N89").Value = -312524292

$ws.Range("H107").Value = 1567.0769
$ws.Range("I107").Value = 1589.2
$ws.Range("J107").Value = 1493.3334
$ws.Range("K107").Value = 1589.2
$ws.Range("L107").Value = 1493.3334
$ws.Range("M107").Value = 330.8
$ws.Range("N107").Value = -5333.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5849904.5
$ws.Range("I31").Value = 1567.9546
$ws.Range("J31").Value = 25644274
$ws.Range("K31").Value = 1567.9546
$ws.Range("L31").Value = 25644274
$ws.Range("M31").Value = -1272.9546
$ws.Range("N31").Value = -25644864

$ws.Range("H34").Value = 5849904.5
$ws.Range("I34").Value = 1567.9546
$ws.Range("J34").Value = 25644274
$ws.Range("K34").Value = 1567.9546
$ws.Range("L34").Value = 25644274
$ws.Range("M34").Value = -1365.9546
$ws.Range("N34").Value = -25644678

$ws.Range("H70").Value = 42567.5
$ws.Range("J70").Value = 42567.5
$ws.Range("L70").Value = 42567.5
$ws.Range("N70").Value = -43197.5

$ws.Range("H73").Value = 42567.5
$ws.Range("J73").Value = 42567.5
$ws.Range("L73").Value = 42567.5
$ws.Range("N73").Value = -44751.5

$ws.Range("H105").Value = 1678.375
$ws.Range("I105").Value = 1678.375
$ws.Range("K105").Value = 1678.375
$ws.Range("M105").Value = 68.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1654.85
$ws.Range("I5").Value = 303
$ws.Range("J5").Value = 4165.4287
$ws.Range("K5").Value = 909
$ws.Range("L5").Value = 12496.2861
$ws.Range("M5").Value = -797
$ws.Range("N5").Value = -12720.2861

$ws.Range("H135").Value = 1654.85
$ws.Range("I135").Value = 303
$ws.Range("J135").Value = 4165.4287
$ws.Range("K135").Value = 2727
$ws.Range("L135").Value = 37488.85830000001
$ws.Range("M135").Value = -192
$ws.Range("N135").Value = -42558.85830000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16280.25
$ws.Range("J15").Value = 16280.25
$ws.Range("L15").Value = 16280.25
$ws.Range("N15").Value = -16856.25

$ws.Range("H70").Value = 9527.886
$ws.Range("I70").Value = 12767.182
$ws.Range("J70").Value = 4046
$ws.Range("K70").Value = 12767.182
$ws.Range("L70").Value = 4046
$ws.Range("M70").Value = -12497.182
$ws.Range("N70").Value = -4586

$ws.Range("H73").Value = 9527.886
$ws.Range("I73").Value = 12767.182
$ws.Range("J73").Value = 4046
$ws.Range("K73").Value = 12767.182
$ws.Range("L73").Value = 4046
$ws.Range("M73").Value = -11831.182
$ws.Range("N73").Value = -5918

$ws.Range("H80").Value = 15876617
$ws.Range("I80").Value = 47621320
$ws.Range("J80").Value = 4264.7144
$ws.Range("K80").Value = 47621320
$ws.Range("L80").Value = 4264.7144
$ws.Range("M80").Value = -47620322
$ws.Range("N80").Value = -6260.7144

$ws.Range("H81").Value = 16280.25
$ws.Range("J81").Value = 16280.25
$ws.Range("L81").Value = 16280.25
$ws.Range("N81").Value = -18276.25

$ws.Range("H83").Value = 15876617
$ws.Range("I83").Value = 47621320
$ws.Range("J83").Value = 4264.7144
$ws.Range("K83").Value = 238106600
$ws.Range("L83").Value = 21323.572
$ws.Range("M83").Value = -238101608
$ws.Range("N83").Value = -31307.572

$ws.Range("H84").Value = 16280.25
$ws.Range("J84").Value = 16280.25
$ws.Range("L84").Value = 48840.75
$ws.Range("N84").Value = -58824.75

$ws.Range("H102").Value = 4891.7334
$ws.Range("I102").Value = 5632.1665
$ws.Range("K102").Value = 5632.1665
$ws.Range("M102").Value = -4010.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 18000
$ws.Range("I40").Value = 26500
$ws.Range("J40").Value = 9500
$ws.Range("K40").Value = 26500
$ws.Range("L40").Value = 9500
$ws.Range("M40").Value = -26364
$ws.Range("N40").Value = -9772

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 2361.1313
$ws.Range("I132").Value = 1721.7161
$ws.Range("J132").Value = 5238.5
$ws.Range("K132").Value = 5165.148300000001
$ws.Range("L132").Value = 15715.5
$ws.Range("M132").Value = -2635.148300000001
$ws.Range("N132").Value = -20775.5

$ws.Range("H136").Value = 13161815
$ws.Range("I136").Value = 16667782
$ws.Range("K136").Value = 50003346
$ws.Range("M136").Value = -50000796

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("J5").Value = 10000000
$ws.Range("L5").Value = 10000000
$ws.Range("N5").Value = -10000224

$ws.Range("H75").Value = 34286
$ws.Range("J75").Value = 34286
$ws.Range("L75").Value = 34286
$ws.Range("N75").Value = -36158

$ws.Range("H78").Value = 34286
$ws.Range("J78").Value = 34286
$ws.Range("L78").Value = 102858
$ws.Range("N78").Value = -112218

$ws.Range("H107").Value = 860.2222
$ws.Range("I107").Value = 1076.0769
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 3228.2307
$ws.Range("L107").Value = 897
$ws.Range("M107").Value = -1308.2307
$ws.Range("N107").Value = -4737

$ws.Range("H132").Value = 1360.8636
$ws.Range("I132").Value = 982.7895
$ws.Range("J132").Value = 3755.3333
$ws.Range("K132").Value = 2948.3685
$ws.Range("L132").Value = 11265.9999
$ws.Range("M132").Value = -418.3685
$ws.Range("N132").Value = -16325.9999

$ws.Range("H135").Value = 43786.25
$ws.Range("J135").Value = 43786.25
$ws.Range("L135").Value = 43786.25
$ws.Range("N135").Value = -53926.25

$ws.Range("H136").Value = 1278.8518
$ws.Range("I136").Value = 979.5217
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2938.5651
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -388.5650999999998
$ws.Range("N136").Value = -14100
